$wb = $excel.ActiveWorkbook

# Both "展览" (Exhibition) and "全部类型" (All Types) sheets carry the same
# data table and need the same updates to column F (想去人数 / "want to go" count).
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 60
    $ws.Range("F3").Value = 407
    $ws.Range("F4").Value = 24
    $ws.Range("F5").Value = 116
}
